# Generate Report for Handoff
# Updates the "latest handoff" timestamps for the b2606e03-... row (row 7)
# across the Overview, zh-cn and de-de sheets to reflect a fresh handoff
# xliff generation run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-29 20:55:28"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-29 20:55:23"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-29 20:55:28"
